$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 (A1, B1) before removing column C
$ws.Cells.Item(1,1).Value = "Скважина+test_data_konkurs (1).xlsx ; УКПГ1АС"
$ws.Cells.Item(1,2).Value = "Газ+test_data_konkurs (1).xlsx ; УКПГ1АС"

# Remove the entire column C (header "Добыча" plus its duplicate data)
$ws.Range("C:C").Delete()
